$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 80; this pushes existing rows 80-87 down to 83-90
$ws.Range("A80:T82").Insert()

# Row 80 (new): Macroferia Regional de Talca / Maule / Caqui / Mankaki, Primera, Region del Maule
$ws.Range("A80").Value = 5
$ws.Range("B80").Value = "Macroferia Regional de Talca"
$ws.Range("C80").Value = "Maule"
$ws.Range("D80").Value = 45106
$ws.Range("E80").Value = 7
$ws.Range("F80").Value = "Fruta"
$ws.Range("G80").Value = 100107
$ws.Range("H80").Value = "Otros"
$ws.Range("I80").Value = 100107001
$ws.Range("J80").Value = "Caqui"
$ws.Range("K80").Value = "Mankaki"
$ws.Range("L80").Value = "Primera"
$ws.Range("M80").Value = 150
$ws.Range("N80").Value = 12000
$ws.Range("O80").Value = 12000
$ws.Range("P80").Value = 12000
$ws.Range("Q80").Value = "$/caja 12 kilos granel"
$ws.Range("R80").Value = "Región del Maule"
$ws.Range("S80").Value = 12000
$ws.Range("T80").Value = 1

# Row 81 (new): Primera, 18 kilos granel, Region del Maule
$ws.Range("A81").Value = 5
$ws.Range("B81").Value = "Macroferia Regional de Talca"
$ws.Range("C81").Value = "Maule"
$ws.Range("D81").Value = 45106
$ws.Range("E81").Value = 7
$ws.Range("F81").Value = "Fruta"
$ws.Range("G81").Value = 100107
$ws.Range("H81").Value = "Otros"
$ws.Range("I81").Value = 100107001
$ws.Range("J81").Value = "Caqui"
$ws.Range("K81").Value = "Mankaki"
$ws.Range("L81").Value = "Primera"
$ws.Range("M81").Value = 120
$ws.Range("N81").Value = 17000
$ws.Range("O81").Value = 17000
$ws.Range("P81").Value = 17000
$ws.Range("Q81").Value = "$/caja 18 kilos granel"
$ws.Range("R81").Value = "Región del Maule"
$ws.Range("S81").Value = 944
$ws.Range("T81").Value = 18

# Row 82 (new): Segunda, 12 kilos granel, Region del Maule
$ws.Range("A82").Value = 5
$ws.Range("B82").Value = "Macroferia Regional de Talca"
$ws.Range("C82").Value = "Maule"
$ws.Range("D82").Value = 45106
$ws.Range("E82").Value = 7
$ws.Range("F82").Value = "Fruta"
$ws.Range("G82").Value = 100107
$ws.Range("H82").Value = "Otros"
$ws.Range("I82").Value = 100107001
$ws.Range("J82").Value = "Caqui"
$ws.Range("K82").Value = "Mankaki"
$ws.Range("L82").Value = "Segunda"
$ws.Range("M82").Value = 100
$ws.Range("N82").Value = 10000
$ws.Range("O82").Value = 10000
$ws.Range("P82").Value = 10000
$ws.Range("Q82").Value = "$/caja 12 kilos granel"
$ws.Range("R82").Value = "Región del Maule"
$ws.Range("S82").Value = 10000
$ws.Range("T82").Value = 1
